$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7730
$ws.Range("I40").Value = 4980
$ws.Range("J40").Value = 8280
$ws.Range("K40").Value = 4980
$ws.Range("L40").Value = 8280
$ws.Range("M40").Value = -4805
$ws.Range("N40").Value = -8630

$ws.Range("H95").Value = 21541.334
$ws.Range("J95").Value = 21541.334
$ws.Range("L95").Value = 21541.334
$ws.Range("N95").Value = -27033.334

$ws.Range("H103").Value = 6511.533
$ws.Range("I103").Value = 788
$ws.Range("J103").Value = 22251.25
$ws.Range("K103").Value = 2364
$ws.Range("L103").Value = 66753.75
$ws.Range("M103").Value = -1778
$ws.Range("N103").Value = -67925.75

$ws.Range("H125").Value = 837.3
$ws.Range("I125").Value = 740.4
$ws.Range("J125").Value = 934.2
$ws.Range("K125").Value = 6663.599999999999
$ws.Range("L125").Value = 8407.800000000001
$ws.Range("M125").Value = -4203.599999999999
$ws.Range("N125").Value = -13327.8

$ws.Range("H137").Value = 713784.4
$ws.Range("I137").Value = 1834985.4
$ws.Range("J137").Value = 2778.8538
$ws.Range("K137").Value = 5504956.199999999
$ws.Range("L137").Value = 8336.561399999999
$ws.Range("M137").Value = -5502406.199999999
$ws.Range("N137").Value = -13436.5614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4713.4575
$ws.Range("I32").Value = 4442.375
$ws.Range("J32").Value = 5896.364
$ws.Range("K32").Value = 4442.375
$ws.Range("L32").Value = 5896.364
$ws.Range("M32").Value = -4155.375
$ws.Range("N32").Value = -6470.364

$ws.Range("H61").Value = 1420.5454
$ws.Range("I61").Value = 1414
$ws.Range("J61").Value = 1450
$ws.Range("K61").Value = 1414
$ws.Range("L61").Value = 1450
$ws.Range("M61").Value = -1202
$ws.Range("N61").Value = -1874

$ws.Range("H136").Value = 1420.5454
$ws.Range("I136").Value = 1414
$ws.Range("J136").Value = 1450
$ws.Range("K136").Value = 4242
$ws.Range("L136").Value = 4350
$ws.Range("M136").Value = -1692
$ws.Range("N136").Value = -9450

$ws.Range("H137").Value = 40261
$ws.Range("J137").Value = 40261
$ws.Range("L137").Value = 40261
$ws.Range("N137").Value = -50461

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 578.9
$ws.Range("I94").Value = 498.42856
$ws.Range("J94").Value = 766.6667
$ws.Range("K94").Value = 498.42856
$ws.Range("L94").Value = 766.6667
$ws.Range("M94").Value = -47.42856
$ws.Range("N94").Value = -1668.6667

$ws.Range("H99").Value = 2426.0625
$ws.Range("I99").Value = 1493
$ws.Range("J99").Value = 3625.7144
$ws.Range("K99").Value = 1493
$ws.Range("L99").Value = 3625.7144
$ws.Range("M99").Value = 5
$ws.Range("N99").Value = -6621.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 262315
$ws.Range("I31").Value = 615149.8
$ws.Range("J31").Value = 3569.4666
$ws.Range("K31").Value = 615149.8
$ws.Range("L31").Value = 3569.4666
$ws.Range("M31").Value = -614854.8
$ws.Range("N31").Value = -4159.4666

$ws.Range("H34").Value = 262315
$ws.Range("I34").Value = 615149.8
$ws.Range("J34").Value = 3569.4666
$ws.Range("K34").Value = 615149.8
$ws.Range("L34").Value = 3569.4666
$ws.Range("M34").Value = -614947.8
$ws.Range("N34").Value = -3973.4666

$ws.Range("H58").Value = 2921.4666
$ws.Range("I58").Value = 1571.4783
$ws.Range("K58").Value = 1571.4783
$ws.Range("M58").Value = -1368.4783

$ws.Range("H86").Value = 2213.36
$ws.Range("I86").Value = 2078.647
$ws.Range("J86").Value = 2499.625
$ws.Range("K86").Value = 2078.647
$ws.Range("L86").Value = 2499.625
$ws.Range("M86").Value = -955.6469999999999
$ws.Range("N86").Value = -4745.625

$ws.Range("H89").Value = 2213.36
$ws.Range("I89").Value = 2078.647
$ws.Range("J89").Value = 2499.625
$ws.Range("K89").Value = 10393.235
$ws.Range("L89").Value = 12498.125
$ws.Range("M89").Value = -4777.235000000001
$ws.Range("N89").Value = -23730.125

$ws.Range("H93").Value = 4990
$ws.Range("I93").Value = 4990
$ws.Range("K93").Value = 4990
$ws.Range("M93").Value = -3118

$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492

$ws.Range("H132").Value = 4084.75
$ws.Range("I132").Value = 2964.5
$ws.Range("J132").Value = 6549.3
$ws.Range("K132").Value = 8893.5
$ws.Range("L132").Value = 19647.9
$ws.Range("M132").Value = -6363.5
$ws.Range("N132").Value = -24707.9

$ws.Range("H136").Value = 2921.4666
$ws.Range("I136").Value = 1571.4783
$ws.Range("K136").Value = 4714.4349
$ws.Range("M136").Value = -2164.4349

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 16602.732
$ws.Range("J34").Value = 11888.272
$ws.Range("L34").Value = 35664.81600000001
$ws.Range("N34").Value = -35832.81600000001

$ws.Range("H48").Value = 3272.7273
$ws.Range("I48").Value = 2500
$ws.Range("J48").Value = 5333.3335
$ws.Range("K48").Value = 7500
$ws.Range("L48").Value = 16000.0005
$ws.Range("M48").Value = -7250
$ws.Range("N48").Value = -16500.0005

$ws.Range("H68").Value = 3045.544
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3045.544
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 9136.632
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -10758.632

$ws.Range("H71").Value = 3045.544
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3045.544
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 27409.896
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -35521.896

$ws.Range("H98").Value = 300
$ws.Range("J98").Value = 300
$ws.Range("L98").Value = 900
$ws.Range("N98").Value = -3896

$ws.Range("H103").Value = 3300
$ws.Range("J103").Value = 3300
$ws.Range("L103").Value = 9900
$ws.Range("N103").Value = -11658

$ws.Range("H107").Value = 13963.104
$ws.Range("J107").Value = 24664.373
$ws.Range("L107").Value = 73993.11900000001
$ws.Range("N107").Value = -77833.11900000001

$ws.Range("H113").Value = 4630445.5
$ws.Range("I113").Value = 667.94116
$ws.Range("K113").Value = 2003.82348
$ws.Range("M113").Value = 166.17652

$ws.Range("H122").Value = 2216.8965
$ws.Range("I122").Value = 730.6316
$ws.Range("J122").Value = 2940.9744
$ws.Range("K122").Value = 6575.6844
$ws.Range("L122").Value = 26468.7696
$ws.Range("M122").Value = -4125.6844
$ws.Range("N122").Value = -31368.7696

$ws.Range("H131").Value = 843.1900000000001
$ws.Range("J131").Value = 855.26044
$ws.Range("L131").Value = 2565.78132
$ws.Range("N131").Value = -12645.78132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6361.756
$ws.Range("I70").Value = 5670.4375
$ws.Range("K70").Value = 5670.4375
$ws.Range("M70").Value = -5400.4375

$ws.Range("H73").Value = 6361.756
$ws.Range("I73").Value = 5670.4375
$ws.Range("K73").Value = 5670.4375
$ws.Range("M73").Value = -4734.4375

$ws.Range("H122").Value = 5113.273
$ws.Range("I122").Value = 3497.2
$ws.Range("J122").Value = 6460
$ws.Range("K122").Value = 10491.6
$ws.Range("L122").Value = 19380
$ws.Range("M122").Value = -8041.599999999999
$ws.Range("N122").Value = -24280

$ws.Range("H132").Value = 3196.5
$ws.Range("I132").Value = 1621.2307
$ws.Range("J132").Value = 4274.316
$ws.Range("K132").Value = 4863.6921
$ws.Range("L132").Value = 12822.948
$ws.Range("M132").Value = -2333.6921
$ws.Range("N132").Value = -17882.948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3528.6086
$ws.Range("I7").Value = 1475.3
$ws.Range("J7").Value = 5108.077
$ws.Range("K7").Value = 1475.3
$ws.Range("L7").Value = 5108.077
$ws.Range("M7").Value = -1363.3
$ws.Range("N7").Value = -5332.077

$ws.Range("H40").Value = 4960.8
$ws.Range("I40").Value = 4765.4287
$ws.Range("J40").Value = 5416.6665
$ws.Range("K40").Value = 4765.4287
$ws.Range("L40").Value = 5416.6665
$ws.Range("M40").Value = -4629.4287
$ws.Range("N40").Value = -5688.6665

$ws.Range("H81").Value = 55316.332
$ws.Range("J81").Value = 55316.332
$ws.Range("L81").Value = 55316.332
$ws.Range("N81").Value = -57312.332

$ws.Range("H84").Value = 55316.332
$ws.Range("J84").Value = 55316.332
$ws.Range("L84").Value = 165948.996
$ws.Range("N84").Value = -175932.996

$ws.Range("H122").Value = 5706.385
$ws.Range("I122").Value = 3325.4285
$ws.Range("J122").Value = 8484.166999999999
$ws.Range("K122").Value = 9976.2855
$ws.Range("L122").Value = 25452.501
$ws.Range("M122").Value = -7526.2855
$ws.Range("N122").Value = -30352.501

$ws.Range("H126").Value = 3528.6086
$ws.Range("I126").Value = 1475.3
$ws.Range("J126").Value = 5108.077
$ws.Range("K126").Value = 4425.9
$ws.Range("L126").Value = 15324.231
$ws.Range("M126").Value = -1955.9
$ws.Range("N126").Value = -20264.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2528.238
$ws.Range("I81").Value = 1840.5
$ws.Range("J81").Value = 3903.7144
$ws.Range("K81").Value = 3681
$ws.Range("L81").Value = 7807.4288
$ws.Range("M81").Value = -2620
$ws.Range("N81").Value = -9929.4288

$ws.Range("H84").Value = 2528.238
$ws.Range("I84").Value = 1840.5
$ws.Range("J84").Value = 3903.7144
$ws.Range("K84").Value = 18405
$ws.Range("L84").Value = 39037.144
$ws.Range("M84").Value = -13101
$ws.Range("N84").Value = -49645.144

$ws.Range("H97").Value = 34940
$ws.Range("J97").Value = 34940
$ws.Range("L97").Value = 34940
$ws.Range("N97").Value = -36922

$ws.Range("H126").Value = 764085.3
$ws.Range("I126").Value = 3252
$ws.Range("J126").Value = 890890.8
$ws.Range("K126").Value = 9756
$ws.Range("L126").Value = 2672672.4
$ws.Range("M126").Value = -7286
$ws.Range("N126").Value = -2677612.4
